$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Data rows 2-50: Wins=70, Losses=92, Ties=0
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 70
    $ws.Cells.Item($r, 31).Value = 92
    $ws.Cells.Item($r, 32).Value = 0
}
